$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '231.44'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.29'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.504'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05552'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.391'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.485'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.141'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.7884'

$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1391'
$ws.Range("E10").Value = '9WazirXWRX'

$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07461'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'

$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03135'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02908'
$ws.Range("E13").Value = '12BitrueCoinBTR'

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09264'
$ws.Range("E14").Value = '13BitMartTokenBMX'

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001663'
$ws.Range("E15").Value = '14BitForexTokenBF'

$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.264'
$ws.Range("E16").Value = '15MCDexMCB'

$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04740'
$ws.Range("E17").Value = '16CoinExTokenCET'

$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005901'
$ws.Range("E18").Value = '17OneONE'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006258'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.005250'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001066'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0001505'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.181'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0005893'
$ws.Range("E27").Value = '26UpBotsUBXTBestin24h'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04029'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007181'

$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1030'
$ws.Range("E42").Value = '41BKEXTokenBKK'

$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003296'
$ws.Range("E43").Value = '42CEJICEJI'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009374'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00005527'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.6776'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.09407'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002108'
